# aggiornamento fino a 02/05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44313, 6, 43, 125.1236687423616),
    @(44314, 4, 43, 125.1236687423616),
    @(44315, 8, 42, 122.2138159809114),
    @(44316, 8, 44, 128.0335215038119),
    @(44317, 7, 44, 128.0335215038119),
    @(44318, 4, 45, 130.9433742652622)
)

$startRow = 239
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $vals[0]
    # carry forward the date-style formatting (border/alignment/numFmt)
    # used by the rest of column A
    $ws.Cells.Item($lastExistingRow, 1).Copy()
    $cellA.PasteSpecial(-4122)

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
